# Adds the new daily record (2020-05-09 / serial 43960) to the
# "Condicion_Pacientes" table, growing it (and the sheet dimension) from
# A1:F40 to A1:F41, and moves the active selection to D43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item("Condicion_Pacientes")

# Grow the table by one row (this also extends the table ref/autofilter
# and the worksheet dimension to A1:F41).
$newRow = $lo.ListRows.Add()

# Copy the date cell's formatting from the row above so the new date
# cell picks up the existing date number format (style index), instead
# of Excel minting a brand-new custom number format.
$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A41").Value = 43960
$ws.Range("B41").Value = 271
$ws.Range("C41").Value = 59
$ws.Range("D41").Value = 237
$ws.Range("E41").Value = 6
$ws.Range("F41").Value = 14

# Move the current selection like the author did in the source commit.
$ws.Range("D43").Select()
